$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update hours worked on the last entry (B33): 2.5 -> 4
$ws.Range("B33").Value = 4

# Move the active selection to B34, mimicking the user entering the value
# in B33 and pressing Enter (cursor advances to the next row).
$ws.Range("B34").Select()

$wb.Save()
